$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update the Runmode for the "Forgot Password" test case from YES to NO
$ws.Range("C2").Value = "NO"

# Move the active selection, mirroring the user's subsequent click
$ws.Activate()
$ws.Range("B8").Select()
